# Fix typo presentation draft: "Xboost" -> "XGBoost"
# Slide 12 ("Models Results") contains a results table (shape "Table 5").
# Row 5, column 1 of that table holds the text "Xboost Regression" split
# across two runs: "Xboost" and " Regression". Only the first run (the
# misspelled product name) needs to be corrected; its formatting
# (including the err="1" spell-check flag) and the trailing " Regression"
# run must be left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table
$cell = $tbl.Cell(5, 1)
$tr = $cell.Shape.TextFrame.TextRange
$tr.Text = "XGBoost"
